$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Date header
Replace-Text "2023-09-08 Friday" "2023-09-09 Saturday"

# Division problems (processed in document order so the "96÷3=" ->
# "15÷5=" / "26÷3=" -> "96÷3=" overlap resolves the same way Word would,
# left-to-right, top-to-bottom, one occurrence at a time).
Replace-Text "11÷3=" "73÷7="
Replace-Text "52÷6=" "23÷7="
Replace-Text "44÷7=" "47÷4="
Replace-Text "19÷2=" "13÷5="
Replace-Text "89÷9=" "14÷5="

Replace-Text "61÷4=" "82÷2="
Replace-Text "66÷7=" "72÷3="
Replace-Text "24÷6=" "38÷4="
Replace-Text "63÷7=" "14÷5="
Replace-Text "92÷3=" "14÷2="

Replace-Text "20÷2=" "11÷6="
Replace-Text "43÷4=" "54÷3="
Replace-Text "83÷7=" "90÷9="
Replace-Text "27÷9=" "94÷3="
Replace-Text "67÷6=" "39÷4="

Replace-Text "18÷7=" "32÷8="
Replace-Text "29÷5=" "32÷7="
Replace-Text "87÷8=" "11÷8="
Replace-Text "70÷4=" "12÷7="
Replace-Text "96÷3=" "15÷5="

Replace-Text "78÷4=" "12÷6="
Replace-Text "33÷3=" "10÷2="
Replace-Text "81÷3=" "22÷6="
Replace-Text "26÷3=" "96÷3="
Replace-Text "53÷8=" "53÷4="
